# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
# Inserts a new "MSME Definition" table (rows 19-23) into the Summary sheet,
# pushing the existing "Sector Distribution Details" block (and everything
# below it) down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows starting at row 19 (shifts old row 22.. down to 28..,
# old row 43 down to 49, etc.)
$ws.Range("A19:A24").EntireRow.Insert()

# New "MSME Definition" table header (row 19) - bold, like the other
# table headers on this sheet (e.g. row 9 / row 15).
$ws.Range("B19").Value = "Number of employees"
$ws.Range("C19").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D19").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B19:D19").Font.Bold = $true

# Micro (row 20)
$ws.Range("A20").Value = "Micro"
$ws.Range("B20").Value = "Group I: <15<br/>Group II: <7"
$ws.Range("C20").Value = "'"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "Group I: Som <150,000<br/>Group II: Som <230,000"

# Small (row 21)
$ws.Range("A21").Value = "Small"
$ws.Range("B21").Value = "Group I: 15-50<br/>Group II: 7-50"
$ws.Range("C21").Value = "'"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "Group I: Som 150,000 - 500,000<br/>Group II: Som 230,000 - 500,000"

# Medium (row 22)
$ws.Range("A22").Value = "Medium"
$ws.Range("B22").Value = "Group I: 51-200<br/>Group II: 16-50"
$ws.Range("C22").Value = "'"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "Som 500,000 - 2,000,000"

# Large (row 23)
$ws.Range("A23").Value = "Large"
$ws.Range("B23").Value = "Group I: >200<br/>Group II: >50"
$ws.Range("C23").Value = "'"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "Som > 2,000,000"

# The row insert doesn't auto-shift the existing hyperlink anchored at the
# old location of the "http://www.stat.kg/rus/part/msp.htm" source line
# (was A43, now A49 after the 6-row shift) - remove it and re-add it in
# the new, correct location, restoring the sheet's custom "HyperLink" look
# (blue underline) afterwards.
$ws.Range("A43").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A49"), "http://www.stat.kg/rus/part/msp.htm") | Out-Null
$ws.Range("A49").Font.Underline = $true
$ws.Range("A49").Font.Color = 16711680
